$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2 = @(401,9,48,67,75,45)
    3 = @(701,3,90,45,97,15)
    4 = @(601,9,60,67,60,42)
    5 = @(201,9,30,15,45,30)
    6 = @(1201,2,10,10,10,10)
    7 = @(501,9,52,30,75,45)
    8 = @(1202,2,10,10,10,10)
    9 = @(1203,3,15,15,15,15)
    10 = @(101,9,30,15,60,15)
    11 = @(901,16,15,45,60,60)
    12 = @(902,1,0,0,0,0)
    13 = @(1001,18,30,75,60,72)
    14 = @(301,6,45,30,60,45)
    15 = @(801,3,67,65,52,45)
    16 = @(502,0,4,0,0,0)
    17 = @(3,0,3,3,3,3)
    18 = @(1101,0,15,30,30,0)
    19 = @(802,0,4,5,4,0)
    20 = @(2,0,2,2,2,2)
    21 = @(1,0,2,2,2,2)
    22 = @(602,0,0,4,0,9)
    23 = @(402,0,0,4,0,0)
    24 = @(702,0,0,0,4,0)
    25 = @(1002,0,0,0,0,9)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $vals[$c]
    }
}
